$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellXml($cell, $paraInnerXml) {
    $xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>$paraInnerXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $cell.Range.InsertXML($xml)
}

# 1. "Y" row / Function cell: "LimeLight Automation" - drop the spell-check
#    proofErr markers and merge the two runs into a single clean run.
Set-CellXml $t.Cell(4,2) '<w:r><w:t>LimeLight Automation</w:t></w:r>'

# 2. "2 - UP" row / Function cell: "Hang" loses the _GoBack bookmark that
#    used to sit right after it.
Set-CellXml $t.Cell(7,2) '<w:r><w:t>Hang</w:t></w:r>'

# 3. "2 - RIGHT" row / Function cell: was empty, now reads "Pickup Piston"
#    and carries the _GoBack bookmark that moved from the "Hang" cell.
Set-CellXml $t.Cell(8,2) '<w:r><w:t>Pickup Piston</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
